$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: azure_namespace, repo_to_import, bitbucket_workspace_id, bitbucket_project_key
$azureNamespace = "repo-migartion/git-project"
$workspaceId = "anilsb06"
$projectKey = "REP"
$repos = @("app-n-pak", "ALMAtasks", "CASAshell", "casa-build-utils", "CASAplotserver")
$n = $repos.Length

# Fill column by column so the shared-string table is built up in the same
# order Excel produced it in (sr index, then namespace, then each repo name,
# then workspace id, then project key).
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $azureNamespace
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $repos[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $workspaceId
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $projectKey
}

# Column widths / formatting (column B was widened to fit the longer
# "repo-migartion/git-project" values; C/D/E nudge slightly wider too)
$ws.Columns.Item(2).ColumnWidth = 24.166666666666668
$ws.Columns.Item(3).ColumnWidth = 29.451822916666668
$ws.Columns.Item(4).ColumnWidth = 24.451822916666668
$ws.Columns.Item(5).ColumnWidth = 24.451822916666668

# Selection matching the saved state
$ws.Range("E2:E6").Select()
